$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BUGLIST (right-hand table) ---------------------------------------
# Row 5 (bug #3): new bug entry "Bug quand on meurt on peut pas restart"
$ws.Range("J5").Value2 = "Bug quand on meurt on peut pas restart"

# --- TODOLIST (left-hand table) ---------------------------------------
# Row 8 (task #6, "Gestion du temps"): state moves from "Started - BUG"
# to "Terminé", restyled to match the other green "Terminé" cells
# (G3/G4/G6: italic, RGB 00B050).
$ws.Range("G8").Value2 = "Terminé"
$ws.Range("G8").Font.Color = 5287936
$ws.Range("G8").Font.Italic = $true

# --- Selection ----------------------------------------------------------
$ws.Range("K5").Select()
